$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the "D-Limonene Ethanol Water" sheet to become the new
#        "EtOAc Ethanol Water (2)" sheet (placed right after the source). ---
$src = $wb.Worksheets.Item("D-Limonene Ethanol Water")
$src.Copy($null, $src)
$new = $wb.Worksheets.Item($src.Index + 1)
$new.Name = "EtOAc Ethanol Water (2)"

# --- 2. Row 2: compound name + measured fractions for EtOAc/Ethanol/Water ---
$new.Range("A2").Value = "EtOAc"
$new.Range("D2").Value = [double]"0.78717661999999999"
$new.Range("E2").Value = [double]"0"
$new.Range("F2").Value = [double]"0.21282338000000001"
$new.Range("G2").Value = [double]"1.890729E-2"
$new.Range("H2").Value = [double]"0"
$new.Range("I2").Value = [double]"0.98109270999999998"
$new.Range("V2").Value = "5 0 5"

# --- Row 3 ---
$new.Range("D3").Value = [double]"0.69091230999999997"
$new.Range("E3").Value = [double]"5.9225680000000003E-2"
$new.Range("F3").Value = [double]"0.249862"
$new.Range("G3").Value = [double]"2.4137039999999998E-2"
$new.Range("H3").Value = [double]"2.3274989999999999E-2"
$new.Range("I3").Value = [double]"0.95258796999999995"
$new.Range("V3").Value = "4,5 0,5 5"

# --- Row 4 (also molar masses in A4:C4) ---
$new.Range("A4").Value = [double]"88.11"
$new.Range("B4").Value = [double]"46.067999999999998"
$new.Range("C4").Value = [double]"18.015000000000001"
$new.Range("D4").Value = [double]"0.60191749000000005"
$new.Range("E4").Value = [double]"0.11138708"
$new.Range("F4").Value = [double]"0.28669541999999998"
$new.Range("G4").Value = [double]"3.0553810000000001E-2"
$new.Range("H4").Value = [double]"4.8392440000000002E-2"
$new.Range("I4").Value = [double]"0.92105375"
$new.Range("V4").Value = "4 1 5"

# --- Row 5 ---
$new.Range("D5").Value = [double]"0.52143351999999998"
$new.Range("E5").Value = [double]"0.15639107999999999"
$new.Range("F5").Value = [double]"0.3221754"
$new.Range("G5").Value = [double]"3.8759540000000002E-2"
$new.Range("H5").Value = [double]"7.6009779999999999E-2"
$new.Range("I5").Value = [double]"0.88523068000000005"
$new.Range("V5").Value = "3,5 1,5 5"

# --- Row 6 (also densities in A6:C6) ---
$new.Range("A6").Value = [double]"0.90200000000000002"
$new.Range("B6").Value = [double]"0.78900000000000003"
$new.Range("C6").Value = [double]"0.997"
$new.Range("D6").Value = [double]"0.44754389"
$new.Range("E6").Value = [double]"0.19476635"
$new.Range("F6").Value = [double]"0.35768976000000002"
$new.Range("G6").Value = [double]"4.9130670000000001E-2"
$new.Range("H6").Value = [double]"0.10652979"
$new.Range("I6").Value = [double]"0.84433954"
$new.Range("V6").Value = "2 3 5"

# --- Row 7 ---
$new.Range("D7").Value = [double]"0.34804383999999999"
$new.Range("E7").Value = [double]"0.23034439000000001"
$new.Range("F7").Value = [double]"0.42161176"
$new.Range("G7").Value = [double]"5.8877369999999998E-2"
$new.Range("H7").Value = [double]"0.13808733000000001"
$new.Range("I7").Value = [double]"0.80303530000000001"
$new.Range("V7").Value = "2,5 2,5 5"

# --- 3. Rows 8-15 on the new sheet keep their number-formatted style but have
#        no data yet (only rows 2-7 were filled in for this new solvent set). ---
$new.Range("D8:I15").ClearContents()
$new.Range("J8:U15").ClearContents()

# --- 4. Old "D-Limonene Ethanol Water" sheet loses the tab-selected flag and
#        its prior multi-cell selection collapses back down to P2. ---
$src.Activate()
$src.Range("P2").Select()

# --- 5. New sheet becomes the active / visible tab, selection on V8. ---
$new.Activate()
$new.Range("V8").Select()
